# cv121262a.xlsx — "correção nos dados e inicio da analise PNAD 2009"
#
# The "dados" sheet had two spurious section-header rows (row 5:
# "situação do domicílio" and, after the first deletion shifts things up,
# the former row 8: "grandes regiões e unidades da federação") that carried
# a label in column A but no data in B:I. Removing those two rows lets the
# real data rows slide up into their place, which is exactly what the
# corrected workbook shows (table now ends at row 38 instead of row 40).
# Also the "unnamed: 1_level_1" sub-header in B2 is renamed to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 sub-header: "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"

# Remove the "situação do domicílio" label-only row.
$ws.Rows.Item(5).Delete() | Out-Null

# After the row-5 deletion, the old "grandes regiões e unidades da
# federação" label-only row (previously row 8) is now row 7.
$ws.Rows.Item(7).Delete() | Out-Null

Write-Output "done"
